$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Responsable_dia"

# Apply autofilter over the data range
$ws.Range("A1:Q38").AutoFilter()

# Register the hidden _FilterDatabase defined name (sheet scoped)
$fd = $ws.Names.Add("_xlnm._FilterDatabase", "='Responsable_dia'!`$A`$1:`$Q`$38")
$fd.Visible = $false

# View options: hide gridlines + freeze header row
$excel.ActiveWindow.DisplayGridlines = $false
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
